# Weekly fruit/vegetable price update: insert a new weekly record as
# row 58, pushing the existing rows 58:66 down to 59:67.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 58 (shifts 58:66 -> 59:67,
# and extends the used range to A1:R67).
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with this week's record.
$ws.Range("A58").Value = 9
$ws.Range("B58").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 44476
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = 100112022
$ws.Range("G58").Value = "Arveja Verde"
$ws.Range("H58").Value = "Perfection"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 25
$ws.Range("K58").Value = 25000
$ws.Range("L58").Value = 26000
$ws.Range("M58").Value = 25480
$ws.Range("N58").Value = "`$/malla 25 kilos"
$ws.Range("O58").Value = "Provincia de Huasco"
$ws.Range("P58").Value = 1019
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = "Hortaliza"
